$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (ESO079-003...): clear the "try again?" note in column L
$ws.Range("L2").Value = $null

# Row 9 (NGC1084_GROUP...): add a new note in column L
$ws.Range("L9").Value = "gotta change the track"

# Row 28 (PGC1092512...): add a new note in column L
$ws.Range("L28").Value = "Redo on bigger machine"

# Row 31 (PGC938075...): replace note text in column L
$ws.Range("L31").Value = "redo with less points"

# Row 10 (NGC1121...): add a new note in column L
$ws.Range("L10").Value = "give up? Ask Denis… or also change the track"

# Row 23 (PGC000902...): replace note text in column L
$ws.Range("L23").Value = "redo with var 4 and 1500 points"

# Update the sheet view: scroll position and active selection cell
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("N27").Select()
